$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "91.199.08"
$ws.Range("E2").Value = "  +4.47%  "

# Row 3
$ws.Range("D3").Value = "3.131.86"
$ws.Range("E3").Value = "  +3.48%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "220.21"
$ws.Range("E5").Value = "  +7.19%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "621.85"
$ws.Range("E6").Value = "  +1.83%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.376"
$ws.Range("E7").Value = "  +5.00%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.929"
$ws.Range("E8").Value = "  +15.26%  "

# Row 9
$ws.Range("E9").Value = "  -0.09%  "

# Row 10
$ws.Range("D10").Value = "3.129.45"
$ws.Range("E10").Value = "  +3.62%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.732"
$ws.Range("E11").Value = "  +26.20%  "

# Row 12
$ws.Range("E12").Value = "  +7.27%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000252"
$ws.Range("E13").Value = "  +10.11%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.48"
$ws.Range("E14").Value = "  +11.87%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.41"
$ws.Range("E15").Value = "  +4.57%  "

# Row 16
$ws.Range("D16").Value = "91.010.91"
$ws.Range("E16").Value = "  +4.14%  "

# Row 17
$ws.Range("D17").Value = "3.724.63"
$ws.Range("E17").Value = "  +3.62%  "

# Row 18
$ws.Range("D18").Value = "3.131.15"
$ws.Range("E18").Value = "  +2.40%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.80"
$ws.Range("E19").Value = "  +21.85%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000226"
$ws.Range("E20").Value = "  +14.75%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.16"
$ws.Range("E21").Value = "  +10.58%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "432.56"
$ws.Range("E22").Value = "  +4.84%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.73"
$ws.Range("E23").Value = "  +10.07%  "

# Row 24
$ws.Range("E24").Value = "  +9.04%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.97"
$ws.Range("E25").Value = "  +13.91%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.37"
$ws.Range("E26").Value = "  +10.34%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "83.97"
$ws.Range("E27").Value = "  +4.87%  "

# Row 28
$ws.Range("D28").Value = "3.306.70"
$ws.Range("E28").Value = "  +2.42%  "

# Row 29
$ws.Range("E29").Value = "  +0.01%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.169"
$ws.Range("E30").Value = "  +9.85%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.90"
$ws.Range("E31").Value = "  +12.89%  "

# Row 32
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "528.95"
$ws.Range("E32").Value = "  +7.16%  "

# Row 33
$ws.Range("B33").Value = "dogwifhat"
$ws.Range("C33").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.88"
$ws.Range("E33").Value = "  +14.02%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.887"
$ws.Range("E34").Value = "  -17.27%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.18"
$ws.Range("E35").Value = "  +11.78%  "

# Row 36
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.142"
$ws.Range("E36").Value = "  +10.42%  "

# Row 37
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.31"
$ws.Range("E37").Value = "  +9.02%  "

# Row 38
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.48"
$ws.Range("E38").Value = "  +7.98%  "

# Row 39
$ws.Range("B39").Value = "PancakeSwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.86"
$ws.Range("E39").Value = "  +5.18%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.30"
$ws.Range("E40").Value = "  +0.65%  "

# Row 41
$ws.Range("E41").Value = "  +0.15%  "

# Row 42
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.148"
$ws.Range("E42").Value = "  +12.17%  "

# Row 43
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0775"
$ws.Range("E43").Value = "  +18.79%  "

# Row 44
$ws.Range("E44").Value = "  +0.03%  "

# Row 45
$ws.Range("B45").Value = "PolygonEcosystemToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.379"
$ws.Range("E45").Value = "  +7.93%  "

# Row 46
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.93"
$ws.Range("E46").Value = "  +9.40%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "144.31"
$ws.Range("E47").Value = "  -2.64%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "44.23"
$ws.Range("E48").Value = "  +2.39%  "

# Row 49
$ws.Range("E49").Value = "  +13.66%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000266"
$ws.Range("E50").Value = "  +29.36%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "168.19"
$ws.Range("E51").Value = "  +10.22%  "
